$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132. This shifts existing rows 132..268 down
# to 133..269, preserving all of their values/formatting, and expands the
# used range (dimension) to A1:R269 automatically.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new record's data. All of the
# "constant" columns (A,B,C,E,F,G,H,I,N,O,Q,R) share the same values as every
# other data row in this sheet.
$ws.Cells.Item(132, 1).Value = 5
$ws.Cells.Item(132, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(132, 3).Value = "Maule"
$ws.Cells.Item(132, 4).Value = 44601
$ws.Cells.Item(132, 5).Value = 7
$ws.Cells.Item(132, 6).Value = 100112023
$ws.Cells.Item(132, 7).Value = "Brócoli"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 800
$ws.Cells.Item(132, 12).Value = 800
$ws.Cells.Item(132, 13).Value = 800
$ws.Cells.Item(132, 14).Value = "$/unidad"
$ws.Cells.Item(132, 15).Value = "Región del Maule"
$ws.Cells.Item(132, 16).Value = 800
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = "Hortaliza"
